# Atualiza bases de dados
# Applies updated values to reports/repasse_fes.xlsx (Sheet1)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (FHEMIG)
$ws.Range("G2").Value = 1169734603.76

# Row 4 (SESP)
$ws.Range("J4").Value = 11986761.95627206

# Row 5 (SETOP)
$ws.Range("E5").Value = 12385707.13
$ws.Range("F5").Value = 6314403.09
$ws.Range("G5").Value = 4341437.76
$ws.Range("J5").Value = 37176172.18999993
$ws.Range("O5").Value = 9342176.19

# Row 6 (UNIMONTES)
$ws.Range("G6").Value = 82632501.23

# Row 7 (FUNED)
$ws.Range("G7").Value = 342342408.91
$ws.Range("J7").Value = 396545580.0037492

# Row 8 (HEMOMINAS)
$ws.Range("G8").Value = 233039239.82
$ws.Range("O8").Value = 245610861

# Row 9 (ESP-MG)
$ws.Range("G9").Value = 11240850.39
$ws.Range("J9").Value = 17280384.03723116
$ws.Range("O9").Value = 11283755.47
$ws.Range("P9").Value = 11251650.39

# Row 10 (FAPEMIG)
$ws.Range("O10").Value = 2917300
